$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 17.71117788448332
$ws.Range("C2").Value = 9.661469030192695
$ws.Range("D2").Value = 5.989860096844992
$ws.Range("E2").Value = 10.72848934388817
$ws.Range("G2").Value = 3.64377339384029
$ws.Range("I2").Value = 23.07526773664336
$ws.Range("L2").Value = 10.01890219430437
$ws.Range("O2").Value = 23.98335557137605
$ws.Range("B3").Value = 17.07286101768012
$ws.Range("C3").Value = 9.259337782351659
$ws.Range("D3").Value = 5.872329053787829
$ws.Range("E3").Value = 10.77420638434491
$ws.Range("G3").Value = 3.646372829342612
$ws.Range("I3").Value = 23.22300447214236
$ws.Range("L3").Value = 9.989179164861625
$ws.Range("O3").Value = 24.06633309052217
$ws.Range("B4").Value = 16.67009769285838
$ws.Range("C4").Value = 9.001820436786817
$ws.Range("D4").Value = 5.800802879549401
$ws.Range("E4").Value = 10.80407056797581
$ws.Range("G4").Value = 3.648052754573518
$ws.Range("I4").Value = 23.32017252398105
$ws.Range("L4").Value = 9.972696776110352
$ws.Range("O4").Value = 24.12497903031254
$ws.Range("B5").Value = 16.50348254199347
$ws.Range("C5").Value = 8.894310697397648
$ws.Range("D5").Value = 5.771860426893751
$ws.Range("E5").Value = 10.81669189474239
$ws.Range("G5").Value = 3.648758496418333
$ws.Range("I5").Value = 23.36138854352884
$ws.Range("L5").Value = 9.966429102599742
$ws.Range("O5").Value = 24.15080290365828
$ws.Range("B6").Value = 16.47567352544783
$ws.Range("C6").Value = 8.876306741751888
$ws.Range("D6").Value = 5.767068245224608
$ws.Range("E6").Value = 10.81881493308665
$ws.Range("G6").Value = 3.648876964335483
$ws.Range("I6").Value = 23.36833010340684
$ws.Range("L6").Value = 9.965415604937482
$ws.Range("O6").Value = 24.15520690088421
$ws.Range("B7").Value = 16.66786039494869
$ws.Range("C7").Value = 9.000380786354823
$ws.Range("D7").Value = 5.800411662243688
$ws.Range("E7").Value = 10.80423895526849
$ws.Range("G7").Value = 3.648062186690457
$ws.Range("I7").Value = 23.3207218271938
$ws.Range("L7").Value = 9.972610424275501
$ws.Range("O7").Value = 24.12531951928944
$ws.Range("B8").Value = 17.49347060447833
$ws.Range("C8").Value = 9.52507022743854
$ws.Range("D8").Value = 5.949229493935099
$ws.Range("E8").Value = 10.74388051718105
$ws.Range("G8").Value = 3.644652316640766
$ws.Range("I8").Value = 23.12486449394033
$ws.Range("L8").Value = 10.00828966140699
$ws.Range("O8").Value = 24.01036251249583
$ws.Range("B9").Value = 19.01724152691436
$ws.Range("C9").Value = 10.46618683881747
$ws.Range("D9").Value = 6.244196950641116
$ws.Range("E9").Value = 10.63973302164413
$ws.Range("G9").Value = 3.638627729650044
$ws.Range("I9").Value = 22.79222518462969
$ws.Range("L9").Value = 10.09206491455631
$ws.Range("O9").Value = 23.84645094375662
$ws.Range("B10").Value = 20.06785006296975
$ws.Range("C10").Value = 11.10005955179574
$ws.Range("D10").Value = 6.460307712340104
$ws.Range("E10").Value = 10.5718577999614
$ws.Range("G10").Value = 3.63460059118349
$ws.Range("I10").Value = 22.57947466802527
$ws.Range("L10").Value = 10.16171618618792
$ws.Range("O10").Value = 23.76410125330559
$ws.Range("B11").Value = 20.52896813591924
$ws.Range("C11").Value = 11.37527455712974
$ws.Range("D11").Value = 6.557987826154276
$ws.Range("E11").Value = 10.54285131291695
$ws.Range("G11").Value = 3.632854239816731
$ws.Range("I11").Value = 22.4896229472925
$ws.Range("L11").Value = 10.19508601962532
$ws.Range("O11").Value = 23.73501657131639
$ws.Range("B12").Value = 20.70103406935698
$ws.Range("C12").Value = 11.47755907197065
$ws.Range("D12").Value = 6.594846745146757
$ws.Range("E12").Value = 10.53213596080057
$ws.Range("G12").Value = 3.632205179504199
$ws.Range("I12").Value = 22.45660058095287
$ws.Range("L12").Value = 10.20795764283464
$ws.Range("O12").Value = 23.72521553318078
$ws.Range("B13").Value = 20.66409190682555
$ws.Range("C13").Value = 11.45561689116999
$ws.Range("D13").Value = 6.586914929805656
$ws.Range("E13").Value = 10.53443175228434
$ws.Range("G13").Value = 3.632344422703104
$ws.Range("I13").Value = 22.4636678522028
$ws.Range("L13").Value = 10.20517516534092
$ws.Range("O13").Value = 23.72727230470037
$ws.Range("B14").Value = 20.54317589646254
$ws.Range("C14").Value = 11.38372854295067
$ws.Range("D14").Value = 6.56102306612856
$ws.Range("E14").Value = 10.54196436871043
$ws.Range("G14").Value = 3.632800596202709
$ws.Range("I14").Value = 22.48688604980221
$ws.Range("L14").Value = 10.19614030425024
$ws.Range("O14").Value = 23.73418588553051
$ws.Range("B15").Value = 20.46877569655037
$ws.Range("C15").Value = 11.33944187149084
$ws.Range("D15").Value = 6.545145420861978
$ws.Range("E15").Value = 10.54661331145718
$ws.Range("G15").Value = 3.633081608353573
$ws.Range("I15").Value = 22.50123860816261
$ws.Range("L15").Value = 10.19063660251876
$ws.Range("O15").Value = 23.73857881258067
$ws.Range("B16").Value = 20.03736496184175
$ws.Range("C16").Value = 11.0818052943662
$ws.Range("D16").Value = 6.453907959081977
$ws.Range("E16").Value = 10.57379106215284
$ws.Range("G16").Value = 3.634716436462514
$ws.Range("I16").Value = 22.58548664469919
$ws.Range("L16").Value = 10.15956869666263
$ws.Range("O16").Value = 23.76617127752787
$ws.Range("B17").Value = 19.76830352593285
$ws.Range("C17").Value = 10.9203549753918
$ws.Range("D17").Value = 6.39774640847906
$ws.Range("E17").Value = 10.59094263129646
$ws.Range("G17").Value = 3.635741231176686
$ws.Range("I17").Value = 22.63894917812513
$ws.Range("L17").Value = 10.14093618194216
$ws.Range("O17").Value = 23.78524983504497
$ws.Range("B18").Value = 19.61197111263319
$ws.Range("C18").Value = 10.8262586628969
$ws.Range("D18").Value = 6.365386750228562
$ws.Range("E18").Value = 10.60098380999289
$ws.Range("G18").Value = 3.636338728324107
$ws.Range("I18").Value = 22.67035128264486
$ws.Range("L18").Value = 10.13037831294515
$ws.Range("O18").Value = 23.79701131811753
$ws.Range("B19").Value = 19.55877338700335
$ws.Range("C19").Value = 10.79418874013353
$ws.Range("D19").Value = 6.354421781841227
$ws.Range("E19").Value = 10.60441382278814
$ws.Range("G19").Value = 3.636542417328632
$ws.Range("I19").Value = 22.68109527269633
$ws.Range("L19").Value = 10.1268311302231
$ws.Range("O19").Value = 23.80112863429962
$ws.Range("B20").Value = 19.79710963574201
$ws.Range("C20").Value = 10.93766973377762
$ws.Range("D20").Value = 6.403731098254921
$ws.Range("E20").Value = 10.58909859828285
$ws.Range("G20").Value = 3.635631306088964
$ws.Range("I20").Value = 22.63319049064577
$ws.Range("L20").Value = 10.14290323329966
$ws.Range("O20").Value = 23.78313727347189
$ws.Range("B21").Value = 20.57876202183958
$ws.Range("C21").Value = 11.40489668385776
$ws.Range("D21").Value = 6.568631980747866
$ws.Range("E21").Value = 10.53974456346608
$ws.Range("G21").Value = 3.632666275185899
$ws.Range("I21").Value = 22.48003904478947
$ws.Range("L21").Value = 10.19878773740736
$ws.Range("O21").Value = 23.73212222547436
$ws.Range("B22").Value = 21.0747070975616
$ws.Range("C22").Value = 11.6989713270653
$ws.Range("D22").Value = 6.67562709997281
$ws.Range("E22").Value = 10.50905540888364
$ws.Range("G22").Value = 3.630799799212771
$ws.Range("I22").Value = 22.38579216692609
$ws.Range("L22").Value = 10.23667882811614
$ws.Range("O22").Value = 23.70585275733855
$ws.Range("B23").Value = 20.81141552028959
$ws.Range("C23").Value = 11.54306357095991
$ws.Range("D23").Value = 6.618605509318959
$ws.Range("E23").Value = 10.52529149953148
$ws.Range("G23").Value = 3.631789466414222
$ws.Range("I23").Value = 22.4355565019533
$ws.Range("L23").Value = 10.21633296919895
$ws.Range("O23").Value = 23.71922368716117
$ws.Range("B24").Value = 19.78409150835873
$ws.Range("C24").Value = 10.92984570378361
$ws.Range("D24").Value = 6.401025640109298
$ws.Range("E24").Value = 10.58993172330801
$ws.Range("G24").Value = 3.635680977270975
$ws.Range("I24").Value = 22.63579191974307
$ws.Range("L24").Value = 10.1420134484946
$ws.Range("O24").Value = 23.78408989305776
$ws.Range("B25").Value = 18.6164214725716
$ws.Range("C25").Value = 10.22144788717174
$ws.Range("D25").Value = 6.164329379514704
$ws.Range("E25").Value = 10.66638839184398
$ws.Range("G25").Value = 3.64018712277908
$ws.Range("I25").Value = 22.87667771648175
$ws.Range("L25").Value = 10.06795426002562
$ws.Range("O25").Value = 23.88414595858163
